$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.767124176025391
$ws.Range("B1").Value = 6.563809871673584
$ws.Range("C1").Value = 5.56760835647583
$ws.Range("D1").Value = 6.659129619598389
$ws.Range("E1").Value = 3.788486242294312
